# "Add files via upload" — the requirements.txt contents were rewritten as a
# flat list of package/module names, one per paragraph, each with the plain
# "eastAsia hint + en-US lang" run formatting, dropping the old Courier-New /
# shaded "numpy<br/>pandas" / "tkinter" styling and adding time, base64,
# image and default_data_transformer. The trailing _GoBack bookmark now
# trails the final (default_data_transformer) paragraph.
#
# Doing this as a single WordOpenXML-fragment swap via Range.InsertXML is far
# more reliable than chaining Find/Replace + InsertParagraphAfter calls to
# hand-edit run/paragraph formatting, so we rebuild the whole body in one
# shot and let Word re-attach the existing sectPr.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ReqParagraph([string]$text, [bool]$pPrHasFont = $true, [bool]$runHasFont = $true) {
    $pPrRpr = if ($pPrHasFont) {
        '<w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/>'
    } else {
        '<w:lang w:val="en-US"/>'
    }
    $runRpr = if ($runHasFont) {
        '<w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/>'
    } else {
        '<w:lang w:val="en-US"/>'
    }
    return "<w:p $wNs><w:pPr><w:rPr>$pPrRpr</w:rPr></w:pPr><w:r><w:rPr>$runRpr</w:rPr><w:t>$text</w:t></w:r></w:p>"
}

$names = @("numpy", "pandas", "time", "tkinter", "base64", "image")

$xml = New-ReqParagraph "streamlit" $true $false
foreach ($name in $names) {
    $xml += New-ReqParagraph $name $true $true
}

# Last paragraph: pPr/rPr keeps only lang (no rFonts hint), and it carries
# the _GoBack bookmark that used to sit at the end of the old "tkinter"
# paragraph.
$xml += "<w:p $wNs><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:lang w:val=`"en-US`"/></w:rPr>" +
        "<w:t>default_data_transformer</w:t></w:r>" +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Replace the whole body story in one go; InsertXML on the full Content
# range leaves sectPr alone and rebuilds just the paragraph content.
$d.Content.InsertXML($xml)

Write-Host ("Paragraphs now: " + $d.Paragraphs.Count)
